$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): add P1=14, Q1=15, matching the style of O1 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)

# --- Data rows (2-25): swap columns I<->K and M<->O, then add P=2, Q=2 ---
for ($r = 2; $r -le 25; $r++) {
    $iVal = $ws.Range("I$r").Value()
    $kVal = $ws.Range("K$r").Value()
    $ws.Range("I$r").Value = $kVal
    $ws.Range("K$r").Value = $iVal

    $mVal = $ws.Range("M$r").Value()
    $oVal = $ws.Range("O$r").Value()
    $ws.Range("M$r").Value = $oVal
    $ws.Range("O$r").Value = $mVal

    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
